$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Shared-string text edit: "Handoff transform failed" -> "Ready for
#    handoff". This string is used by Overview!B2/C2 and by the "Status"
#    column (B2) on both the zh-cn and de-de sheets, so update all of them
#    so they end up pointing at the same (deduplicated) new string.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")
$wsZh.Range("B2").Value = "Ready for handoff"
$wsDe.Range("B2").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: record a completed handoff for row 2.
#    C2 = Latest Handoff File (new hyperlinked cell)
#    D2 = Latest Handoff Datetime
#    H2 = Handoff Reason: "Ignored" -> "Include"
# ---------------------------------------------------------------------------
$wsZh.Range("D2").Value = "2016-01-18 03:12:40"
$wsZh.Range("H2").Value = "Include"

# Re-insert the A3 hyperlink after adding the new C2 hyperlink so the
# <hyperlinks> order becomes A2, C2, A3 (matching a freshly appended
# relationship for C2 while A3's relationship moves to the end).
foreach ($hlink in $wsZh.Hyperlinks) {
    $addr = $hlink.Range.Address()
    if ($addr -eq '$A$3') {
        $hlink.Delete()
    }
}
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/654eba25ce0d35b586304d619e0dc4f228817e27/e2e/c0152aff-a42c-467a-88ab-51f83d00bffb.f7ca64998f331c82bdc2197ad84b0ea535dd57fb.zh-cn.xlf", "", "", "c0152aff-a42c-467a-88ab-51f83d00bffb.f7ca64998f331c82bdc2197ad84b0ea535dd57fb.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/654eba25ce0d35b586304d619e0dc4f228817e27/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# 3) de-de sheet: same shape of edit, different file/datetime.
# ---------------------------------------------------------------------------
$wsDe.Range("D2").Value = "2016-01-18 03:12:51"
$wsDe.Range("H2").Value = "Include"

foreach ($hlink in $wsDe.Hyperlinks) {
    $addr = $hlink.Range.Address()
    if ($addr -eq '$A$3') {
        $hlink.Delete()
    }
}
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/654eba25ce0d35b586304d619e0dc4f228817e27/e2e/c0152aff-a42c-467a-88ab-51f83d00bffb.f7ca64998f331c82bdc2197ad84b0ea535dd57fb.de-de.xlf", "", "", "c0152aff-a42c-467a-88ab-51f83d00bffb.f7ca64998f331c82bdc2197ad84b0ea535dd57fb.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/654eba25ce0d35b586304d619e0dc4f228817e27/.localization-config", "", "", ".localization-config")

Write-Output "applied localization handoff report updates"
